$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, pushing existing rows 15-40 down to 16-41.
$ws.Rows.Item(15).Insert()

# The newly inserted row 15 is blank; seed it with a copy of the row that is
# now at row 16 (the original row 15 content), then overwrite the handful
# of cells that differ for the new record.
$ws.Range("A16:T16").Copy($ws.Range("A15:T15"))

$ws.Cells.Item(15, 4).Value = [DateTime]"2023-12-12"   # D15 (Fecha) -> serial 45272
$ws.Cells.Item(15, 13).Value = 90                        # M15 (Volumen)
$ws.Cells.Item(15, 14).Value = 22000                      # N15 (Precio minimo)
$ws.Cells.Item(15, 15).Value = 22000                      # O15 (Precio maximo)
$ws.Cells.Item(15, 16).Value = 22000                      # P15 (Precio promedio ponderado)
$ws.Cells.Item(15, 19).Value = 2200                       # S15 (Precio $/Kg)
